# Auto-generated: refresh cached market-price / profit figures on multiple
# sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), per scheduled runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1217328.9
$ws.Range("I6").Value = 1217328.9
$ws.Range("K6").Value = 3651986.7
$ws.Range("M6").Value = -3651874.7
$ws.Range("H9").Value = 195.64285
$ws.Range("I9").Value = 206.36363
$ws.Range("J9").Value = 156.33333
$ws.Range("K9").Value = 206.36363
$ws.Range("L9").Value = 156.33333
$ws.Range("M9").Value = -37.36363
$ws.Range("N9").Value = -494.33333
$ws.Range("H10").Value = 10210
$ws.Range("I10").Value = 420
$ws.Range("J10").Value = 20000
$ws.Range("K10").Value = 420
$ws.Range("L10").Value = 20000
$ws.Range("M10").Value = -127
$ws.Range("N10").Value = -20586
$ws.Range("H12").Value = 791.6667
$ws.Range("I12").Value = 700
$ws.Range("K12").Value = 700
$ws.Range("M12").Value = -530
$ws.Range("H29").Value = 5403.8887
$ws.Range("J29").Value = 6500.6665
$ws.Range("L29").Value = 19501.9995
$ws.Range("N29").Value = -20063.9995
$ws.Range("H54").Value = 510038
$ws.Range("I54").Value = 510038
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 510038
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -509552
$ws.Range("N54").ClearContents()
$ws.Range("H55").Value = 541.5
$ws.Range("I55").Value = 287.14285
$ws.Range("J55").Value = 795.8570999999999
$ws.Range("K55").Value = 287.14285
$ws.Range("L55").Value = 795.8570999999999
$ws.Range("M55").Value = -73.14285000000001
$ws.Range("N55").Value = -1223.8571
$ws.Range("H62").Value = 7287.2
$ws.Range("I62").Value = 5051.1665
$ws.Range("J62").Value = 8777.888999999999
$ws.Range("K62").Value = 5051.1665
$ws.Range("L62").Value = 8777.888999999999
$ws.Range("M62").Value = -4427.1665
$ws.Range("N62").Value = -10025.889
$ws.Range("H65").Value = 7287.2
$ws.Range("I65").Value = 5051.1665
$ws.Range("J65").Value = 8777.888999999999
$ws.Range("K65").Value = 25255.8325
$ws.Range("L65").Value = 43889.44499999999
$ws.Range("M65").Value = -22135.8325
$ws.Range("N65").Value = -50129.44499999999
$ws.Range("H80").Value = 1870.2
$ws.Range("I80").Value = 410
$ws.Range("J80").Value = 2032.4445
$ws.Range("K80").Value = 1230
$ws.Range("L80").Value = 6097.333500000001
$ws.Range("M80").Value = -232
$ws.Range("N80").Value = -8093.333500000001
$ws.Range("H83").Value = 1870.2
$ws.Range("I83").Value = 410
$ws.Range("J83").Value = 2032.4445
$ws.Range("K83").Value = 3690
$ws.Range("L83").Value = 18292.0005
$ws.Range("M83").Value = 1302
$ws.Range("N83").Value = -28276.0005
$ws.Range("H106").Value = 18825.5
$ws.Range("I106").Value = 10867.637
$ws.Range("K106").Value = 10867.637
$ws.Range("M106").Value = -10236.637
$ws.Range("H107").Value = 1968.1364
$ws.Range("I107").Value = 1870.2307
$ws.Range("J107").Value = 2109.5557
$ws.Range("K107").Value = 1870.2307
$ws.Range("L107").Value = 2109.5557
$ws.Range("M107").Value = 49.76929999999993
$ws.Range("N107").Value = -5949.5557
$ws.Range("H125").Value = 1029.8572
$ws.Range("I125").Value = 1032
$ws.Range("J125").Value = 1029.5
$ws.Range("K125").Value = 9288
$ws.Range("L125").Value = 9265.5
$ws.Range("M125").Value = -6828
$ws.Range("N125").Value = -14185.5
$ws.Range("H138").Value = 4768.5713
$ws.Range("I138").Value = 3902.6667
$ws.Range("J138").Value = 5220.3477
$ws.Range("K138").Value = 11708.0001
$ws.Range("L138").Value = 15661.0431
$ws.Range("M138").Value = -6568.000100000001
$ws.Range("N138").Value = -25941.0431
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1282.2391
$ws.Range("I32").Value = 1283.7046
$ws.Range("K32").Value = 1283.7046
$ws.Range("M32").Value = -996.7046
$ws.Range("H36").Value = 9967.286
$ws.Range("I36").Value = 5954.3
$ws.Range("J36").Value = 19999.75
$ws.Range("K36").Value = 5954.3
$ws.Range("L36").Value = 19999.75
$ws.Range("M36").Value = -5608.3
$ws.Range("N36").Value = -20691.75
$ws.Range("H45").Value = 76927640
$ws.Range("I45").Value = 111114050
$ws.Range("J45").Value = 8210.25
$ws.Range("K45").Value = 111114050
$ws.Range("L45").Value = 8210.25
$ws.Range("M45").Value = -111113673
$ws.Range("N45").Value = -8964.25
$ws.Range("H63").Value = 4380.92
$ws.Range("I63").Value = 2313.923
$ws.Range("J63").Value = 6620.1665
$ws.Range("K63").Value = 2313.923
$ws.Range("L63").Value = 6620.1665
$ws.Range("M63").Value = -1627.923
$ws.Range("N63").Value = -7992.1665
$ws.Range("H66").Value = 4380.92
$ws.Range("I66").Value = 2313.923
$ws.Range("J66").Value = 6620.1665
$ws.Range("K66").Value = 11569.615
$ws.Range("L66").Value = 33100.8325
$ws.Range("M66").Value = -8137.614999999998
$ws.Range("N66").Value = -39964.8325
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2384.5
$ws.Range("I107").Value = 2146.7273
$ws.Range("J107").Value = 5000
$ws.Range("K107").Value = 2146.7273
$ws.Range("L107").Value = 5000
$ws.Range("M107").Value = -226.7273
$ws.Range("N107").Value = -8840
$ws.Range("H134").Value = 1936.881
$ws.Range("I134").Value = 1135.0294
$ws.Range("J134").Value = 5344.75
$ws.Range("K134").Value = 3405.0882
$ws.Range("L134").Value = 16034.25
$ws.Range("M134").Value = -870.0881999999997
$ws.Range("N134").Value = -21104.25
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1531.8334
$ws.Range("I6").Value = 1531.8334
$ws.Range("K6").Value = 1531.8334
$ws.Range("M6").Value = -1418.8334
$ws.Range("H31").Value = 21189.639
$ws.Range("I31").Value = 2362.6667
$ws.Range("J31").Value = 36140.47
$ws.Range("K31").Value = 2362.6667
$ws.Range("L31").Value = 36140.47
$ws.Range("M31").Value = -2067.6667
$ws.Range("N31").Value = -36730.47
$ws.Range("H34").Value = 21189.639
$ws.Range("I34").Value = 2362.6667
$ws.Range("J34").Value = 36140.47
$ws.Range("K34").Value = 2362.6667
$ws.Range("L34").Value = 36140.47
$ws.Range("M34").Value = -2160.6667
$ws.Range("N34").Value = -36544.47
$ws.Range("H86").Value = 6995.1
$ws.Range("I86").Value = 3499.8333
$ws.Range("K86").Value = 3499.8333
$ws.Range("M86").Value = -2376.8333
$ws.Range("H89").Value = 6995.1
$ws.Range("I89").Value = 3499.8333
$ws.Range("K89").Value = 17499.1665
$ws.Range("M89").Value = -11883.1665
$ws.Range("H93").Value = 13499.5
$ws.Range("I93").Value = 13499.5
$ws.Range("K93").Value = 13499.5
$ws.Range("M93").Value = -11627.5
$ws.Range("H99").Value = 3502.75
$ws.Range("I99").Value = 3349.5
$ws.Range("J99").Value = 3656
$ws.Range("K99").Value = 3349.5
$ws.Range("L99").Value = 3656
$ws.Range("M99").Value = -1851.5
$ws.Range("N99").Value = -6652
$ws.Range("H105").Value = 3272.7778
$ws.Range("I105").Value = 5399.8
$ws.Range("J105").Value = 2454.6924
$ws.Range("K105").Value = 5399.8
$ws.Range("L105").Value = 2454.6924
$ws.Range("M105").Value = -3652.8
$ws.Range("N105").Value = -5948.6924
$ws.Range("H126").Value = 3502.75
$ws.Range("I126").Value = 3349.5
$ws.Range("J126").Value = 3656
$ws.Range("K126").Value = 10048.5
$ws.Range("L126").Value = 10968
$ws.Range("M126").Value = -7578.5
$ws.Range("N126").Value = -15908
$ws.Range("H132").Value = 2595.7317
$ws.Range("I132").Value = 2058.1843
$ws.Range("K132").Value = 6174.5529
$ws.Range("M132").Value = -3644.5529
$ws.Range("H141").Value = 65000
$ws.Range("J141").Value = 65000
$ws.Range("L141").Value = 65000
$ws.Range("N141").Value = -75360
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 57227.43
$ws.Range("I2").Value = 41.357143
$ws.Range("J2").Value = 114413.5
$ws.Range("K2").Value = 248.142858
$ws.Range("L2").Value = 686481
$ws.Range("M2").Value = -135.142858
$ws.Range("N2").Value = -686707
$ws.Range("H47").Value = 141665.67
$ws.Range("I47").Value = 87500
$ws.Range("K47").Value = 262500
$ws.Range("M47").Value = -262069
$ws.Range("H93").Value = 15954.5
$ws.Range("I93").Value = 10350
$ws.Range("J93").Value = 18756.75
$ws.Range("K93").Value = 31050
$ws.Range("L93").Value = 56270.25
$ws.Range("M93").Value = -29178
$ws.Range("N93").Value = -60014.25
$ws.Range("H132").Value = 7107.7144
$ws.Range("I132").Value = 6000
$ws.Range("K132").Value = 54000
$ws.Range("M132").Value = -51470
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 70.75
$ws.Range("I2").Value = 67.55
$ws.Range("J2").Value = 78.75
$ws.Range("K2").Value = 67.55
$ws.Range("L2").Value = 78.75
$ws.Range("M2").Value = 45.45
$ws.Range("N2").Value = -304.75
$ws.Range("H132").Value = 3781.9644
$ws.Range("I132").Value = 2612
$ws.Range("J132").Value = 18991.5
$ws.Range("K132").Value = 7836
$ws.Range("L132").Value = 56974.5
$ws.Range("M132").Value = -5306
$ws.Range("N132").Value = -62034.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1463.8
$ws.Range("J16").Value = 2575
$ws.Range("L16").Value = 2575
$ws.Range("N16").Value = -2915
$ws.Range("H132").Value = 4546.853
$ws.Range("I132").Value = 3228.2727
$ws.Range("J132").Value = 6964.25
$ws.Range("K132").Value = 9684.8181
$ws.Range("L132").Value = 20892.75
$ws.Range("M132").Value = -7154.8181
$ws.Range("N132").Value = -25952.75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5979.9375
$ws.Range("I81").Value = 4169.643
$ws.Range("K81").Value = 8339.286
$ws.Range("M81").Value = -7278.286
$ws.Range("H84").Value = 5979.9375
$ws.Range("I84").Value = 4169.643
$ws.Range("K84").Value = 41696.43
$ws.Range("M84").Value = -36392.43
$ws.Range("H126").Value = 3668.8845
$ws.Range("I126").Value = 2285.6924
$ws.Range("J126").Value = 5052.077
$ws.Range("K126").Value = 6857.0772
$ws.Range("L126").Value = 15156.231
$ws.Range("M126").Value = -4387.0772
$ws.Range("N126").Value = -20096.231
$ws.Range("H132").Value = 5215.4043
$ws.Range("I132").Value = 2619.7693
$ws.Range("J132").Value = 17869.125
$ws.Range("K132").Value = 7859.3079
$ws.Range("L132").Value = 53607.375
$ws.Range("M132").Value = -5329.3079
$ws.Range("N132").Value = -58667.375
$ws.Range("H136").Value = 4742.222
$ws.Range("I136").Value = 4079.1785
$ws.Range("K136").Value = 12237.5355
$ws.Range("M136").Value = -9687.5355
